$wb = $excel.ActiveWorkbook

# "binek" sheet - sheetId 1 / rId1 / sheet1.xml
$ws = $wb.Worksheets.Item("binek")

# Update variable label in A7
$ws.Range("A7").Value = "kredi_kullanan_arac_orani"

# Update values
$ws.Range("B2").Value = 0.7
$ws.Range("B6").Value = 120000
$ws.Range("B7").Value = 0.52

# Update the active selection on the sheet
$ws.Range("B7").Select()
